$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '69.659.46'
$ws.Range("E2").Value = '  +3.02%  '
$ws.Range("D3").Value = '3.393.18'
$ws.Range("E3").Value = '  +4.41%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '191.73'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +4.12%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '594.14'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.41%  '
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("E8").Value = '  +0.19%  '
$ws.Range("E9").Value = '  +2.71%  '
$ws.Range("E11").Value = '  +1.97%  '
$ws.Range("D12").Value = '3.984.01'
$ws.Range("E12").Value = '  +4.85%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.136'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.75%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '28.81'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +4.09%  '
$ws.Range("D15").Value = '69.607.08'
$ws.Range("E15").Value = '  +2.94%  '
$ws.Range("E16").Value = '  +1.70%  '
$ws.Range("D17").Value = '3.391.52'
$ws.Range("E17").Value = '  +5.55%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '451.49'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +14.28%  '
$ws.Range("E19").Value = '  +1.58%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.85'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +2.32%  '
$ws.Range("E21").Value = '  +3.55%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '76.39'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +6.91%  '
$ws.Range("E23").Value = '  +0.20%  '
$ws.Range("E24").Value = '  +1.41%  '
$ws.Range("E25").Value = '  +4.30%  '
$ws.Range("E26").Value = '  +2.00%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.51'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.64%  '
$ws.Range("E28").Value = '  -0.16%  '
$ws.Range("E29").Value = '  +2.84%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '23.53'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +3.91%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '5.67'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +2.04%  '
$ws.Range("E32").Value = '  +3.05%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '7.01'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.03%  '
$ws.Range("E34").Value = '  +0.02%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.58'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +7.00%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '165.43'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +2.51%  '
$ws.Range("E37").Value = '  +2.95%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '28.40'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +6.51%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.818'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.31%  '
$ws.Range("E40").Value = '  +1.83%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.61'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +2.09%  '
$ws.Range("D42").Value = '2.755.50'
$ws.Range("E42").Value = '  +5.31%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.53'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +2.09%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '25.67'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +3.84%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '41.17'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.30%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '342.20'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +2.22%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0286'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +2.58%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '33.15'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +7.81%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.02'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +5.73%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '6.35'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.44%  '
